# The original single run "Fadsfsafd" must become three separate runs:
#   "Fadsfs", "XXXXX", "afd"
# A plain Find/Replace (or Range.InsertAfter/InsertBefore) would still end
# up as a single run, because adjacent runs that share identical run
# formatting get coalesced back together when the document is saved.
# To guarantee three distinct <w:r> elements (with no extra run
# properties), locate the run's range and overwrite it with explicit
# WordprocessingML via Range.InsertXML, which inserts exactly the runs
# we specify.

$d = $word.ActiveDocument

$r = $d.Content
$found = $r.Find.Execute("Fadsfsafd", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)

$pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
       '<pkg:xmlData>' + `
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
       '<w:body><w:p>' + `
       '<w:r><w:t>Fadsfs</w:t></w:r>' + `
       '<w:r><w:t>XXXXX</w:t></w:r>' + `
       '<w:r><w:t>afd</w:t></w:r>' + `
       '</w:p></w:body></w:document>' + `
       '</pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($pkg)
